# Updated Cursos Settings MonthlyRaport added for printing
$d = $word.ActiveDocument

# --- Title paragraph: #[ewidencja-title] font shrinks from 12pt (sz 24) to 6pt (sz 12) ---
$title = $d.Paragraphs.Item(1)
$title.Range.Font.Size = 6
$title.Range.Font.SizeBi = 6

# --- First empty paragraph right after the table: sz 18 -> sz 16 (9pt -> 8pt) ---
# NOTE: this must happen *before* any Table/Tables access below - touching the
# Tables collection first throws off later Paragraphs.Item(N) lookups in this
# runtime. The paragraph itself has no run, and Range.Font setters only
# rewrite the paragraph-mark (pPr/rPr) formatting when the range actually
# contains a run, so temporarily add then remove a placeholder character.
$tailPara = $d.Paragraphs.Item(28)
$tailPara.Range.InsertAfter("x")
$tailPara.Range.Font.Size = 8
$tailPara.Range.Font.SizeBi = 8
$tailRange = $d.Paragraphs.Item(28).Range
$removeRange = $d.Range($tailRange.Start, $tailRange.End - 1)
$removeRange.Delete()

# --- Table: uniform column widths, row heights, and smaller table font ---
$t = $d.Tables.Item(1)

# Shrink every run (header + data rows) from 9pt (sz 18) to 6pt (sz 12).
# Using the whole table Range also rewrites the paragraph-mark run
# properties (pPr/rPr) for every cell paragraph, matching the diff.
$t.Range.Font.Size = 6
$t.Range.Font.SizeBi = 6

# Make every column the same width: 1528 twips = 76.4 pt.
for ($c = 1; $c -le $t.Columns.Count; $c++) {
    $t.Columns.Item($c).Width = 76.4
}

# Row heights: header row 795 twips = 39.75 pt, data row 1071 twips = 53.55 pt.
$t.Rows.Item(1).Height = 39.75
$t.Rows.Item(2).Height = 53.55
